# Update the "تاريخ تفعيل الخدمة" (service activation date) entry in D4:
# the stray "01/28/ 2021" text value is removed, leaving the cell blank.
# (Once no cell references that shared string any more, the workbook
# automatically drops it from sharedStrings.xml on save.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("D4")
$target.ClearContents()

# Reset the cell back to the default/general style (instead of the
# inherited date-format style) by pulling the plain formatting from an
# already-unstyled cell, so no new style entries get added to the sheet.
$blankStyleSource = $ws.Range("G1")
$blankStyleSource.Copy()
$target.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The workbook was last saved with the C8 cell selected/active.
[void]$ws.Range("C8").Select()
